# Update header text on the existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add a new "PO Forecast" sheet after the existing sheets
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row, matching the bold / centered / top-aligned / bordered style
# used by the other sheets' header rows.
$hdr = $wsForecast.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$poForecastData = @(
    @(44934.99999999999, 14, -34.42054019009382, 60.63784897538051),
    @(44962.99999999999, 16, -34.20108106586881, 65.93316836903486),
    @(44990.99999999999, 19, -27.79347065480449, 68.57863083156219),
    @(45004.99999999999, 20, -29.02868305338998, 67.15787096834961),
    @(45011.99999999999, 20, -26.1947197363981, 68.14175615923571),
    @(45025.99999999999, 21, -26.41816898181564, 67.92889152350139),
    @(45032.99999999999, 22, -27.02246299611982, 68.52525121014499),
    @(45039.99999999999, 23, -25.45304357094516, 73.82377219212019),
    @(45046.99999999999, 23, -26.9216321537827, 73.81314782526992),
    @(45060.99999999999, 24, -23.25683910347569, 71.30597178968971),
    @(45074.99999999999, 25, -19.17607314932688, 73.28065985723978),
    @(45081.99999999999, 26, -24.63579889278532, 74.77273443343007),
    @(45088.99999999999, 26, -20.90167181528979, 72.18491542982349),
    @(45095.99999999999, 27, -17.39775521193259, 73.8106411695363),
    @(45116.99999999999, 29, -17.86193538001364, 78.43853028313323),
    @(45123.99999999999, 29, -18.71507796523728, 79.91598684151835),
    @(45130.99999999999, 30, -18.49808211726618, 79.39773272712426),
    @(45137.99999999999, 30, -20.78156777712043, 75.92078491248292),
    @(45144.99999999999, 31, -19.61519559947045, 80.99886296131294),
    @(45151.99999999999, 32, -14.32885630365605, 80.73624171616461),
    @(45186.99999999999, 34, -13.58116781013862, 84.24403054809908),
    @(45193.99999999999, 35, -14.9296593856765, 84.70429518625676),
    @(45200.99999999999, 35, -12.16058556851759, 86.63774814158914),
    @(45221.99999999999, 37, -9.927705606499197, 86.38352232116038),
    @(45228.99999999999, 38, -9.422573498424448, 86.98871895167629),
    @(45235.99999999999, 38, -9.360796368355565, 87.0565981696847),
    @(45249.99999999999, 39, -12.23477099506716, 87.16095576016718),
    @(45256.99999999999, 40, -9.832888099556566, 86.03110278126199),
    @(45263.99999999999, 41, -7.399798204074009, 89.61976440403539),
    @(45270.99999999999, 41, -5.157495179352763, 86.93991289635436),
    @(45277.99999999999, 42, -3.932206616699959, 92.97819874056019),
    @(45305.99999999999, 44, -3.09143027555462, 90.93232336859592),
    @(45312.99999999999, 44, -1.350288082451456, 93.06342700508232),
    @(45319.99999999999, 45, -2.798788158078372, 91.6770411984111),
    @(45333.99999999999, 46, -3.565927848353067, 97.05842795931635),
    @(45340.99999999999, 47, -1.5214046217052, 94.29279284715409),
    @(45382.99999999999, 50, 4.603945563236996, 102.1770114146954),
    @(45396.99999999999, 51, 3.766475300784025, 99.40124568873249),
    @(45410.99999999999, 52, 4.508290590528943, 99.44782288474715),
    @(45438.99999999999, 55, 6.548321835870046, 102.8211260625629),
    @(45452.99999999999, 56, 10.7056086689372, 107.7787692879117),
    @(45459.99999999999, 56, 7.169584108710771, 105.6972637236469),
    @(45466.99999999999, 57, 8.469754751379639, 103.5795260275413),
    @(45571.99999999999, 65, 15.67108368682052, 112.8485834675167),
    @(45578.99999999999, 66, 20.20062408227799, 111.5846218148192),
    @(45599.99999999999, 68, 16.8398579395725, 118.3489770417834),
    @(45641.99999999999, 71, 25.46959490495845, 119.93120154061),
    @(45648.99999999999, 71, 23.7592121092132, 120.0396597963157),
    @(45655.99999999999, 72, 24.76412840606424, 119.4060966969063),
    @(45662.99999999999, 73, 28.17293644303616, 122.7143375696434),
    @(45669.99999999999, 73, 30.45454487063269, 119.856729178247),
    @(45676.99999999999, 74, 26.27993907317831, 121.0093245887259),
    @(45683.99999999999, 74, 23.34957262822785, 124.3595774306044),
    @(45690.99999999999, 75, 27.61932138362594, 122.0189752018011),
    @(45697.99999999999, 75, 27.81579348453374, 120.0765760384221)
)

$rowCount = $poForecastData.Length
for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $poForecastData[$i]
    $r = $i + 2
    $dateCell = $wsForecast.Cells.Item($r, 1)
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $dateCell.Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
}

Write-Output "Added PO Forecast sheet with $rowCount data rows"

# Restore the originally active sheet/tab selection.
$wsWeekly.Activate()
